# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" sheet between "2021-Q4" and "总计" with the
#    quarter's fund-holding detail rows (same shape as the "2021-Q4" sheet).
# 2) Prepend a "2022-Q1" row to the "总计" roll-up sheet (existing
#    "2021-Q4" row shifts down one row).
#
# Formatting note: this workbook's header row + index column (col A) use a
# bold/bordered/centered style (xf index 2 in the original file). COM has no
# "apply style index N" verb, so we copy it over with Copy/PasteSpecial
# (xlPasteFormats = -4122) from cells that already carry it. Also: assigning
# a numeric-looking string straight to .Value (e.g. "004616" or "94.50")
# gets auto-coerced to a number and loses the leading/trailing zeros, so
# those are entered with a leading "'" (force-text) and then have their
# format re-pasted from a plain/unstyled cell, which clears the resulting
# quote-prefix marker again and leaves a clean text cell.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

$ref = $wb.Worksheets.Item("2021-Q4")   # reference sheet: already-correct styles live here
$refHeaderCell = $ref.Range("B1")       # style 2: bold + border + centered
$refIndexCell  = $ref.Range("A2")       # style 2: bold + border + centered
$refPlainCell  = $ref.Range("C2")       # style 0: default/no override

# ---------------------------------------------------------------------
# 1) New sheet "2022-Q1", placed right after "2021-Q4"
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($null, $ref)
$q1.Name = "2022-Q1"

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data rows
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'004616"
$q1.Range("C2").Value = "中欧电子信息产业沪港深股票A"
$q1.Range("D2").Value = "'14.54"
$q1.Range("E2").Value = "'92.26"
$q1.Range("F2").Value = "'3.08"
$q1.Range("G2").Value = "'0.4478"
$q1.Range("H2").Value = 9

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'005763"
$q1.Range("C3").Value = "中欧电子信息产业沪港深股票C"
$q1.Range("D3").Value = "'7.73"
$q1.Range("E3").Value = "'92.26"
$q1.Range("F3").Value = "'3.08"
$q1.Range("G3").Value = "'0.2381"
$q1.Range("H3").Value = 9

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "'501219"
$q1.Range("C4").Value = "华夏智胜先锋股票（LOF）A"
$q1.Range("D4").Value = "'3.61"
$q1.Range("E4").Value = "'94.50"
$q1.Range("F4").Value = "'1.13"
$q1.Range("G4").Value = "'0.0408"
$q1.Range("H4").Value = 3

$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "'014198"
$q1.Range("C5").Value = "华夏智胜先锋股票（LOF）C"
$q1.Range("D5").Value = "'1.30"
$q1.Range("E5").Value = "'94.50"
$q1.Range("F5").Value = "'1.13"
$q1.Range("G5").Value = "'0.0147"
$q1.Range("H5").Value = 3

# Re-apply the real formats, which also scrubs the quote-prefix left by the
# force-text entry above.
$refHeaderCell.Copy()
$q1.Range("B1:H1").PasteSpecial($xlPasteFormats)

$refIndexCell.Copy()
$q1.Range("A2:A5").PasteSpecial($xlPasteFormats)

$refPlainCell.Copy()
$q1.Range("B2:H5").PasteSpecial($xlPasteFormats)
$q1.Range("H2:H5").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 2) "总计" roll-up sheet: shift the old 2021-Q4 total down to row 3 and
#    write the new 2022-Q1 total into row 2.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 1.98

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.74

# A2 already carried style 2 before this script ran; copy it onto the new
# A3 index cell so both rows match.
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial($xlPasteFormats)
